$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46042
$ws.Range("B2").Value = 103.8
$ws.Range("C2").Value = 98.13
$ws.Range("D2").Value = 90.31
$ws.Range("E2").Value = 86.27
$ws.Range("F2").Value = 83.83
$ws.Range("G2").Value = 87.52
$ws.Range("H2").Value = 96.94
$ws.Range("I2").Value = 109.99
$ws.Range("J2").Value = 123.86
$ws.Range("K2").Value = 118.19
$ws.Range("L2").Value = 98.20999999999999
$ws.Range("M2").Value = 94.54000000000001
$ws.Range("N2").Value = 93.56
$ws.Range("O2").Value = 86.15000000000001
$ws.Range("P2").Value = 86.05
$ws.Range("Q2").Value = 86.98999999999999
$ws.Range("R2").Value = 92.12
$ws.Range("S2").Value = 108.99
$ws.Range("T2").Value = 125.1
$ws.Range("U2").Value = 127.25
$ws.Range("V2").Value = 127.01
$ws.Range("W2").Value = 120.83
$ws.Range("X2").Value = 108.33
$ws.Range("Y2").Value = 97.23999999999999
$ws.Range("Z2").Value = 102.13
$ws.Range("AB2").Value = 113.36
$ws.Range("AD2").Value = 126.18
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 123.92
$ws.Range("AG2").Value = "1h-23h"
